$d = $word.ActiveDocument

# --- Change 1: paragraph 1 "EXNO" + ":1" -> single run "EXNO:2 ... Visualizing Time series Data",
#     and move the _GoBack bookmark here (it currently sits at the end of a "figsize=(10, 6))" paragraph).

# Remove the existing _GoBack bookmark from wherever it currently lives.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# Replace the whole first paragraph's content (merges the "EXNO" / ":1" runs and
# drops the proofErr gramStart/gramEnd markers that bracketed them).
$r = $d.Range(0, 7)
$r.Delete()
$r2 = $d.Range(0, 0)
$r2.InsertParagraphBefore()
$p1 = $d.Paragraphs.Item(1)
$newTitle = "EXNO:2                                     Visualizing Time series Data"
# Append a one-character sentinel so the bookmark-insertion point below is not the
# exact end-of-paragraph boundary (which the engine treats specially / ambiguously).
$p1.Range.Text = $newTitle + "X"

$endPos = $p1.Range.End - 2
$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the sentinel character now that the bookmark is anchored right after the text.
$sentinelPos = $p1.Range.End - 2
$d.Range($sentinelPos, $sentinelPos + 1).Delete()

# --- Change 2: merge the "'Date" / "'] = " runs into a single "'Date'] = " run.

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "df['Date'] = pd.to_datetime(df['Date'])`r") {
        $target = $para
        break
    }
}
$s = $target.Range.Start
$mergeRange = $d.Range($s + 3, $s + 13)
# Force an actual text change first (Word skips the edit - and so the run merge -
# if the replacement text is identical to what's already there), then restore the
# real text; the second assignment lands in a single freshly-merged run.
$mergeRange.Text = "XYZPLACEHOLDER"
$mergeRange2 = $d.Range($s + 3, $s + 3 + 14)
$mergeRange2.Text = "'Date'] = "

Write-Output "done"
